$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.565.54"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.275.32"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.26"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.88"
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.17"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.20"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.631.76"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.280.10"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.525.43"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.92"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.06"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.60"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.40"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.92"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.52"
$ws.Range("E36").Value = "  -3.02%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0688"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.001.65"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0277"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.01"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("E47").Value = "  -7.92%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("E49").Value = "  +5.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.47"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.499.89"
$ws.Range("E51").Value = "  -1.01%  "
